$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "file date" column (B) values -------------------------
# Old primary-key dates (various formats of 19/07/2023, 17/06/2023, etc.)
# are replaced with the new primary-key date 01/08/2023 (and its alternate
# textual representations). A leading apostrophe is used so the cells keep
# being stored as text (preserving their existing quoted-text number
# format / style) instead of being reinterpreted as real dates.
$ws.Range("B1").Value  = "'01/08/2023"
$ws.Range("B2").Value  = "'01/08/2023"
$ws.Range("B3").Value  = "'01/08/2023"
$ws.Range("B4").Value  = "'01/08/2023"
$ws.Range("B5").Value  = "'01/08/2023"
$ws.Range("B6").Value  = "'01/08/2023"
$ws.Range("B7").Value  = "'01/08/2023"
$ws.Range("B8").Value  = "'01/08/2023"
$ws.Range("B10").Value = "'01/08/2023"
$ws.Range("B11").Value = "'01/08/2023"
$ws.Range("B12").Value = "'2023/08/01"
$ws.Range("B14").Value = "'08/01/2023"
$ws.Range("B16").Value = "'01/08/2023"
$ws.Range("B17").Value = "'01/08/2023"

# --- Column width changes ----------------------------------------------
# Column C grows, and the previously-merged C5:C6 width definition is
# split so column E gets its own (wider) width while column F keeps the
# original shared width.
$ws.Columns.Item(3).ColumnWidth = 14.5
$ws.Columns.Item(5).ColumnWidth = 17

# --- Selection moves from B6 to B1 --------------------------------------
[void]$ws.Range("B1").Select()
